$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$S0 = "Ementa atual:"
$S1 = "Ementa modificada (dados modificados em vermelho):"
$S2 = "LOM3084"
$S3 = "Nome:"
$S4 = " Inspeção e Ensaios Não Destrutivos"
$S5 = "Name:"
$S6 = "Inspection and Nondestructive Testing"
$S7 = "Créditos-aula:"
$S8 = "2"
$S9 = "Créditos-trabalho"
$S10 = "0"
$S11 = "Carga horária:"
$S12 = "30 h"
$S13 = "Ativação:"
$S14 = "01/01/2020"
$S15 = "Semestre ideal:"
$S16 = "EM-8"
$S17 = "Objetivos:"
$S18 = "Abordar os aspectos fundamentais do assunto, visando a formação de engenheiros habilitados para a escolha adequada do tipo de ensaio não-destrutivo para o controle e detecção de defeitos em estruturas e componentes de engenharia, sua condução e a correta interpretação dos resultados."
$S19 = "Objectives:"
$S20 = "Docentes responsáveis:"
$S21 = "3586455 - Cassius Olivio Figueiredo Terra Ruchert"
$S22 = "5840793 - Sérgio Schneider"
$S23 = "Programa resumido:"
$S24 = "Ensaios Não-Destrutivos (END): princípios e tipos. Nomenclatura. Técnicas especiais de END. Classificação, normas e especificações."
$S25 = "Short syllabus:"
$S26 = "Programa:"
$S27 = "Introdução aos Ensaios Não-Destrutivos (END). Tipos mais comuns de END. Líquidos penetrantes: princípios, materiais, procedimentos, aplicações e limitações, normas e especificações, critérios de aceitação, qualificações. Ensaio por ultrassom: fundamentos, equipamentos, aplicações e limitações, transdutores, critérios de aceitação, normas e especificações, qualificações. Ensaio por partículas magnéticas: propriedades magnéticas, fundamentos do ensaio, equipamentos, aplicações e limitações, critérios de aceitação, normas e especificações, qualificação. Ensaio por correntes parasitas: fundamentos, equipamentos, aplicações e limitações, critérios de aceitação, normas e especificações, qualificação. Ensaio radiográfico: radiações, princípio, fontes de radiação, Proteção radiológica (Normas CNEN), interpretação, descontinuidades típicas, aplicações e limitações, critérios de aceitação, qualificação. Métodos especiais de END: Emissão de Barkhausen, ensaio de vibração."
$S28 = "Syllabus:"
$S29 = "Avaliação:"
$S30 = "Método:"
$S31 = "Serão aplicadas duas avaliações: P1 e P2"
$S32 = "Critério:"
$S33 = "Conceito Final = (P1+P2)/2 ou (P1+MS)/2       (MS = média de seminários)"
$S34 = "Norma de recuperação:"
$S35 = "A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação"
$S36 = "Bibliografia:"
$S37 = "API Standards. American Petroleum Institute; 2011. 
ASM Handbook, Vol. 17: Non destructive evaluation and quality control. American Society for Materials; 1989. 
ASME Handbook, Vol. : Non destructive evaluation and quality control. American Society for Materials; 1989. 
ASNT Handbook, Vol. 10: Nondestructive Testing Overview. American Society for Nondestructive Testing, 1993. 
CARTZ, L. Nondestructive Testing. American Society for Testing and Materials, 1995."
$S38 = "Requisitos:"
$S39 = "LOM3107 -  Mecânica dos Sólidos Deformáveis  (Requisito fraco)
"
$S40 = "LOM3109 -  Fundamentos da Física Moderna  (Requisito fraco)
"

$ws.Range("B10").Value = $S18
$ws.Range("C10").Value = $S18
$ws.Range("A13").ClearContents()
$ws.Range("B13").Value = $S21
$ws.Range("C13").Value = $S21
$ws.Range("A14").ClearContents()
$ws.Range("B14").Value = $S22
$ws.Range("C14").Value = $S22
$ws.Range("A15").Value = $S23
$ws.Range("B15").Value = $S24
$ws.Range("C15").Value = $S24
$ws.Range("A16").Value = $S25
$ws.Range("A17").Value = $S26
$ws.Range("B17").Value = $S27
$ws.Range("C17").Value = $S27
$ws.Range("A18").Value = $S28
$ws.Range("B18").ClearContents()
$ws.Range("C18").ClearContents()
$ws.Range("A19").Value = $S29
$ws.Range("B19").ClearContents()
$ws.Range("C19").ClearContents()
$ws.Range("A20").Value = $S30
$ws.Range("B20").Value = $S31
$ws.Range("C20").Value = $S31
$ws.Range("A21").Value = $S32
$ws.Range("B21").Value = $S33
$ws.Range("C21").Value = $S33
$ws.Range("A22").Value = $S34
$ws.Range("B22").Value = $S35
$ws.Range("C22").Value = $S35
$ws.Range("A23").Value = $S36
$ws.Range("B23").Value = $S37
$ws.Range("C23").Value = $S37
$ws.Range("A24").Value = $S38
$ws.Range("B24").ClearContents()
$ws.Range("C24").ClearContents()
$ws.Range("B25").Value = $S39
$ws.Range("C25").Value = $S39
$ws.Range("B26").Value = $S40
$ws.Range("C26").Value = $S40

$ws.Rows.Item(13).AutoFit()
$ws.Rows.Item(14).AutoFit()
$ws.Rows.Item(15).RowHeight = 60
$ws.Rows.Item(16).RowHeight = 60
$ws.Rows.Item(17).RowHeight = 120
$ws.Rows.Item(18).RowHeight = 120
$ws.Rows.Item(19).AutoFit()
$ws.Rows.Item(21).RowHeight = 60
$ws.Rows.Item(22).RowHeight = 60
$ws.Rows.Item(23).RowHeight = 120
$ws.Rows.Item(24).AutoFit()
$ws.Rows.Item(25).RowHeight = 30
$ws.Rows.Item(26).RowHeight = 30
